$d = $word.ActiveDocument

$d.Content.Find.Execute("84$([char]0x00D7)74=6216", $true, $false, $false, $false, $false, $true, 1, $false, "79$([char]0x00D7)73=5767", 2) | Out-Null
$d.Content.Find.Execute("51$([char]0x00D7)16=816", $true, $false, $false, $false, $false, $true, 1, $false, "17$([char]0x00D7)70=1190", 2) | Out-Null
$d.Content.Find.Execute("88$([char]0x00D7)52=4576", $true, $false, $false, $false, $false, $true, 1, $false, "92$([char]0x00D7)17=1564", 2) | Out-Null
$d.Content.Find.Execute("88$([char]0x00D7)85=7480", $true, $false, $false, $false, $false, $true, 1, $false, "32$([char]0x00D7)86=2752", 2) | Out-Null
$d.Content.Find.Execute("35$([char]0x00D7)66=2310", $true, $false, $false, $false, $false, $true, 1, $false, "73$([char]0x00D7)65=4745", 2) | Out-Null
$d.Content.Find.Execute("84$([char]0x00D7)78=6552", $true, $false, $false, $false, $false, $true, 1, $false, "16$([char]0x00D7)89=1424", 2) | Out-Null
$d.Content.Find.Execute("13$([char]0x00D7)36=468", $true, $false, $false, $false, $false, $true, 1, $false, "52$([char]0x00D7)71=3692", 2) | Out-Null
$d.Content.Find.Execute("38$([char]0x00D7)60=2280", $true, $false, $false, $false, $false, $true, 1, $false, "31$([char]0x00D7)92=2852", 2) | Out-Null
$d.Content.Find.Execute("23$([char]0x00D7)11=253", $true, $false, $false, $false, $false, $true, 1, $false, "70$([char]0x00D7)38=2660", 2) | Out-Null
$d.Content.Find.Execute("69$([char]0x00D7)39=2691", $true, $false, $false, $false, $false, $true, 1, $false, "97$([char]0x00D7)92=8924", 2) | Out-Null
$d.Content.Find.Execute("70$([char]0x00D7)73=5110", $true, $false, $false, $false, $false, $true, 1, $false, "82$([char]0x00D7)48=3936", 2) | Out-Null
$d.Content.Find.Execute("13$([char]0x00D7)40=520", $true, $false, $false, $false, $false, $true, 1, $false, "42$([char]0x00D7)50=2100", 2) | Out-Null
$d.Content.Find.Execute("13$([char]0x00D7)60=780", $true, $false, $false, $false, $false, $true, 1, $false, "11$([char]0x00D7)62=682", 2) | Out-Null
$d.Content.Find.Execute("75$([char]0x00D7)99=7425", $true, $false, $false, $false, $false, $true, 1, $false, "13$([char]0x00D7)32=416", 2) | Out-Null
$d.Content.Find.Execute("42$([char]0x00D7)67=2814", $true, $false, $false, $false, $false, $true, 1, $false, "78$([char]0x00D7)76=5928", 2) | Out-Null
$d.Content.Find.Execute("15$([char]0x00D7)55=825", $true, $false, $false, $false, $false, $true, 1, $false, "84$([char]0x00D7)54=4536", 2) | Out-Null
$d.Content.Find.Execute("23$([char]0x00D7)93=2139", $true, $false, $false, $false, $false, $true, 1, $false, "25$([char]0x00D7)42=1050", 2) | Out-Null
$d.Content.Find.Execute("95$([char]0x00D7)49=4655", $true, $false, $false, $false, $false, $true, 1, $false, "79$([char]0x00D7)49=3871", 2) | Out-Null
$d.Content.Find.Execute("35$([char]0x00D7)50=1750", $true, $false, $false, $false, $false, $true, 1, $false, "21$([char]0x00D7)94=1974", 2) | Out-Null
$d.Content.Find.Execute("95$([char]0x00D7)82=7790", $true, $false, $false, $false, $false, $true, 1, $false, "74$([char]0x00D7)20=1480", 2) | Out-Null
$d.Content.Find.Execute("75$([char]0x00D7)38=2850", $true, $false, $false, $false, $false, $true, 1, $false, "14$([char]0x00D7)98=1372", 2) | Out-Null
$d.Content.Find.Execute("12$([char]0x00D7)41=492", $true, $false, $false, $false, $false, $true, 1, $false, "65$([char]0x00D7)49=3185", 2) | Out-Null
$d.Content.Find.Execute("76$([char]0x00D7)92=6992", $true, $false, $false, $false, $false, $true, 1, $false, "17$([char]0x00D7)27=459", 2) | Out-Null
$d.Content.Find.Execute("92$([char]0x00D7)85=7820", $true, $false, $false, $false, $false, $true, 1, $false, "92$([char]0x00D7)43=3956", 2) | Out-Null
$d.Content.Find.Execute("90$([char]0x00D7)96=8640", $true, $false, $false, $false, $false, $true, 1, $false, "94$([char]0x00D7)51=4794", 2) | Out-Null
